# Update customer report to v22
# Remove the "Marketplace ID", "Marketplace Name" and "Environment" columns
# from the customer-list report (columns I:J and C in the original layout),
# shifting the remaining columns left, and refresh the AutoFilter /
# _FilterDatabase range to match the new (narrower) data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "Marketplace ID" (I) and "Marketplace Name" (J) first, then
# "Environment" (C) -- deleting the higher-indexed columns first keeps the
# lower column letters stable while we work.
$ws.Range("I1:J1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()

# The data now spans A1:T1 (20 columns). Re-apply the AutoFilter over the
# new extent (turn it off first since the stored filter range otherwise
# keeps pointing at the old, wider range).
$ws.AutoFilterMode = $false
$ws.Range("A1:S1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# shrunk filter range as well.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$A`$1:`$S`$1"
    }
}
